$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 198, pushing the existing rows 198:226 down to 200:228
# (same as selecting rows 198:199 and doing Insert > Shift cells down).
$ws.Rows("198:199").Insert()

# New row 198: Zapallo italiano, Primera, week of 2021-11-22
$ws.Cells.Item(198, 1).Value = 1
$ws.Cells.Item(198, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(198, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(198, 4).Value = 44522
$ws.Cells.Item(198, 5).Value = 15
$ws.Cells.Item(198, 6).Value = 100112032
$ws.Cells.Item(198, 7).Value = "Zapallo italiano"
$ws.Cells.Item(198, 8).Value = "Huracán"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 120
$ws.Cells.Item(198, 11).Value = 4000
$ws.Cells.Item(198, 12).Value = 4500
$ws.Cells.Item(198, 13).Value = 4250
$ws.Cells.Item(198, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(198, 16).Value = 61
$ws.Cells.Item(198, 17).Value = 70
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# New row 199: Zapallo italiano, Segunda, week of 2021-11-22
$ws.Cells.Item(199, 1).Value = 1
$ws.Cells.Item(199, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(199, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(199, 4).Value = 44522
$ws.Cells.Item(199, 5).Value = 15
$ws.Cells.Item(199, 6).Value = 100112032
$ws.Cells.Item(199, 7).Value = "Zapallo italiano"
$ws.Cells.Item(199, 8).Value = "Huracán"
$ws.Cells.Item(199, 9).Value = "Segunda"
$ws.Cells.Item(199, 10).Value = 140
$ws.Cells.Item(199, 11).Value = 3500
$ws.Cells.Item(199, 12).Value = 4000
$ws.Cells.Item(199, 13).Value = 3750
$ws.Cells.Item(199, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(199, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(199, 16).Value = 38
$ws.Cells.Item(199, 17).Value = 100
$ws.Cells.Item(199, 18).Value = "Hortaliza"
